$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 42 (Caso 5883 / CONGRESO AV. 2699) — all subsequent rows shift up by one.
$ws.Rows.Item(42).Delete()
